$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154, shifting rows 154:200 down to 155:201
$ws.Rows("154").Insert()

# Populate the new row 154 with the new record's data
$ws.Range("A154").Value = 5
$ws.Range("B154").Value = "Macroferia Regional de Talca"
$ws.Range("C154").Value = "Maule"
$ws.Range("D154").Value = 44627
$ws.Range("E154").Value = 7
$ws.Range("F154").Value = 100112045
$ws.Range("G154").Value = "Zapallo"
$ws.Range("H154").Value = "Camote"
$ws.Range("I154").Value = "1a nueva(o)"
$ws.Range("J154").Value = 900
$ws.Range("K154").Value = 300
$ws.Range("L154").Value = 300
$ws.Range("M154").Value = 300
$ws.Range("N154").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O154").Value = "Región del Maule"
$ws.Range("P154").Value = 300
$ws.Range("Q154").Value = 1
$ws.Range("R154").Value = "Hortaliza"
